$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) cells whose new value looks like a plain number. ---
# Excel would silently coerce these to numeric cells (losing formatting like
# trailing zeros or literal multi-dot grouping), so force a Text number format
# on just these cells before writing the literal string.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.56"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.00"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.46"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.94"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.759"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.73"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.73"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.75"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.90"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.35"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.62"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0624"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0987"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.81"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0214"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.58"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.41"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.01"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.49"

# --- Price (column D) cells whose new value is not a plain number ---
# (e.g. multi-dot grouped values or subscripted text) stay text natively.
$ws.Range("D2").Value = "37.716.21"
$ws.Range("D3").Value = "2.077.73"
$ws.Range("D12").Value = "2.372.11"
$ws.Range("D17").Value = "2.069.78"
$ws.Range("D18").Value = "37.662.04"
$ws.Range("D21").Value = "0.0₃0819"
$ws.Range("D45").Value = "1.454.37"

# --- Volume/percentage change (column E) cells ---
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("E27").Value = "  +10.40%  "
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("E37").Value = "  +4.71%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("E40").Value = "  +5.74%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("E51").Value = "  +7.60%  "
